$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44487
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 30
$ws.Range("N2").Value = 23000
$ws.Range("O2").Value = 24000
$ws.Range("P2").Value = 23500
$ws.Range("S2").Value = 2350
$ws.Range("D3").Value = 44839
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 120
$ws.Range("N3").Value = 25000
$ws.Range("O3").Value = 26000
$ws.Range("P3").Value = 25500
$ws.Range("S3").Value = 2550
$ws.Range("D4").Value = 44452
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 21000
$ws.Range("O4").Value = 22000
$ws.Range("P4").Value = 21500
$ws.Range("S4").Value = 2150
$ws.Range("D5").Value = 44446
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 21000
$ws.Range("O5").Value = 22000
$ws.Range("P5").Value = 21500
$ws.Range("S5").Value = 2150
$ws.Range("D6").Value = 44461
$ws.Range("L6").Value = "Especial"
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 31000
$ws.Range("O6").Value = 32000
$ws.Range("P6").Value = 31500
$ws.Range("S6").Value = 3150
$ws.Range("D7").Value = 44461
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 30
$ws.Range("N7").Value = 30000
$ws.Range("O7").Value = 30000
$ws.Range("P7").Value = 30000
$ws.Range("S7").Value = 3000
$ws.Range("D8").Value = 44874
$ws.Range("L8").Value = "Especial"
$ws.Range("M8").Value = 30
$ws.Range("N8").Value = 25000
$ws.Range("O8").Value = 25000
$ws.Range("P8").Value = 25000
$ws.Range("S8").Value = 2500
$ws.Range("D9").Value = 44874
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 80
$ws.Range("N9").Value = 23000
$ws.Range("O9").Value = 24000
$ws.Range("P9").Value = 23500
$ws.Range("S9").Value = 2350
$ws.Range("D10").Value = 44841
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 60
$ws.Range("N10").Value = 23000
$ws.Range("O10").Value = 24000
$ws.Range("P10").Value = 23500
$ws.Range("S10").Value = 2350
$ws.Range("D11").Value = 44868
$ws.Range("L11").Value = "Especial"
$ws.Range("M11").Value = 60
$ws.Range("N11").Value = 26000
$ws.Range("O11").Value = 26000
$ws.Range("P11").Value = 26000
$ws.Range("S11").Value = 2600
$ws.Range("D12").Value = 44460
$ws.Range("L12").Value = "Especial"
$ws.Range("M12").Value = 60
$ws.Range("N12").Value = 31000
$ws.Range("O12").Value = 32000
$ws.Range("P12").Value = 31500
$ws.Range("S12").Value = 3150
$ws.Range("D13").Value = 44460
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 30
$ws.Range("N13").Value = 30000
$ws.Range("O13").Value = 30000
$ws.Range("P13").Value = 30000
$ws.Range("S13").Value = 3000
$ws.Range("D14").Value = 44448
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 60
$ws.Range("N14").Value = 21000
$ws.Range("O14").Value = 22000
$ws.Range("P14").Value = 21500
$ws.Range("S14").Value = 2150
$ws.Range("D15").Value = 44447
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 60
$ws.Range("N15").Value = 21000
$ws.Range("O15").Value = 22000
$ws.Range("P15").Value = 21500
$ws.Range("S15").Value = 2150
$ws.Range("D16").Value = 44848
$ws.Range("L16").Value = "Especial"
$ws.Range("M16").Value = 60
$ws.Range("N16").Value = 24000
$ws.Range("O16").Value = 25000
$ws.Range("P16").Value = 24500
$ws.Range("S16").Value = 2450
$ws.Range("D17").Value = 44848
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 120
$ws.Range("N17").Value = 21000
$ws.Range("O17").Value = 22000
$ws.Range("P17").Value = 21500
$ws.Range("S17").Value = 2150
